$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = 30
$ws.Range("B9").Value = "Edit1"
$ws.Range("C9").Value = "riya-morankar"
$ws.Range("D9").Value = "N/A"
$ws.Range("E9").Value = "edit1 to main"

# Force the date-looking value to be stored as literal text (matches the
# source row's "2025-06-17" style values, which are plain strings, not
# real dates) instead of being auto-parsed into a date serial number.
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "2025-06-18"
$ws.Range("F9").Style = "Normal"
